# "Played With Laser Calculations"
# Update the Laser sheet: change Lens Output Half Angle and Altitude inputs,
# bump the packet size, and introduce an explicit "Bits per byte" row (10b/8b
# line-encoding overhead) feeding the bad-packet-rate formulas instead of the
# previously hard-coded *8 multiplier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Laser")

# Insert a new row above the old "Bad Packet Rate (RZ-OOK)(Night)" row (was
# row 25) to hold the new "Bits per byte" input. This shifts everything from
# the old row 25 down by one (rows 25-30 -> 26-31), and Excel auto-adjusts
# the relative formulas referencing rows above the insertion point.
$ws.Rows.Item(25).Insert()

# New row 25: Bits per byte input (10 bits per byte due to 10b/8b encoding)
$ws.Range("A25").Value = "Bits per byte"
$ws.Range("B25").Value = 10
$ws.Range("C25").Value = "10b/8b encoding"
$ws.Range("B25").Style = "Input"
$ws.Range("B25").NumberFormat = "0 ""b"""

# Update the (now shifted) bad-packet-rate formulas to use B25 instead of a
# hard-coded 8.
$ws.Range("B26").Formula = "=1-(1-B22)^(B24*B25)"
$ws.Range("B27").Formula = "=1-(1-B23)^(B24*B25)"

# Core input changes on the Laser sheet.
$ws.Range("B6").Value = 0.5     # Lens Output Half Angle: 1 -> 0.5
$ws.Range("B7").Value = 600     # Altitude: 400 -> 600
$ws.Range("B24").Value = 1024   # Packet Size: 1000 -> 1024

# Update the view selection to match the edited area.
$ws.Range("B7").Select()

$wb.Save()
